# Bump,Bump,Bump.xlsx -> lyrics swapped to "The Streak" (Ray Stevens).
# Sheet shrinks from A1:B114 to A1:B75 (39 trailing rows removed) and
# rows 2-75 get new Section/Content text per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded trailing rows (old Chorus 15-21 etc.) first so the
# used range / dimension collapses to A1:B75 once the new content is written.
$ws.Range("A76:B114").EntireRow.Delete()

$ws.Cells.Item(2, 1).Value = 'Verse 1'
$ws.Cells.Item(2, 2).Value = '[]'
$ws.Cells.Item(3, 1).Value = 'Verse 2'
$ws.Cells.Item(3, 2).Value = '[''Hello, everyone, this is your action news reporter with all the news that is news across the nation, on the scene at the supermarket.'']'
$ws.Cells.Item(4, 1).Value = 'Verse 3'
$ws.Cells.Item(4, 2).Value = '[''Yeah, I did.'']'
$ws.Cells.Item(5, 1).Value = 'Verse 4'
$ws.Cells.Item(5, 2).Value = '["I''s standin'' overe there by the tomaters, and here he come, running through the pole beans, through the fruits and vegetables, nekkid as a jay bird."]'
$ws.Cells.Item(6, 1).Value = 'Verse 5'
$ws.Cells.Item(6, 2).Value = '[''And I hollered over t\'' Ethel, I said, "Don\''t look, Ethel!"'']'
$ws.Cells.Item(7, 1).Value = 'Verse 6'
$ws.Cells.Item(7, 2).Value = '["But it''s too late, she''d already been incensed."]'
$ws.Cells.Item(8, 1).Value = 'Verse 7'
$ws.Cells.Item(8, 2).Value = '[''Here he comes, look at that, look at that'']'
$ws.Cells.Item(9, 1).Value = 'Verse 8'
$ws.Cells.Item(9, 2).Value = '[''There he goes, look at that, look at that'']'
$ws.Cells.Item(10, 1).Value = 'Verse 9'
$ws.Cells.Item(10, 2).Value = '["And he ain''t wearin'' no clothes"]'
$ws.Cells.Item(11, 1).Value = 'Verse 10'
$ws.Cells.Item(11, 2).Value = '[''(Look at that, look at that)'']'
$ws.Cells.Item(12, 1).Value = 'Verse 11'
$ws.Cells.Item(12, 2).Value = '[''Of his anatomy'']'
$ws.Cells.Item(13, 1).Value = 'Verse 12'
$ws.Cells.Item(13, 2).Value = '["Invitin'' public critique"]'
$ws.Cells.Item(14, 1).Value = 'Verse 13'
$ws.Cells.Item(14, 2).Value = '[''Yeah, I did.'']'
$ws.Cells.Item(15, 1).Value = 'Verse 14'
$ws.Cells.Item(15, 2).Value = '["I''s just in here gettin my car checked, he just appeared out of the traffic."]'
$ws.Cells.Item(16, 1).Value = 'Verse 15'
$ws.Cells.Item(16, 2).Value = '["Come streakin'' around the grease rack there, didn''t have nothin'' on but a smile."]'
$ws.Cells.Item(17, 1).Value = 'Verse 16'
$ws.Cells.Item(17, 2).Value = '["I looked in there, and Ethel was gettin'' her a cold drink."]'
$ws.Cells.Item(18, 1).Value = 'Verse 17'
$ws.Cells.Item(18, 2).Value = '[]'
$ws.Cells.Item(19, 1).Value = 'Verse 18'
$ws.Cells.Item(19, 2).Value = '["She''d already been mooned."]'
$ws.Cells.Item(20, 1).Value = 'Verse 19'
$ws.Cells.Item(20, 2).Value = '[''Flashed her right there in front of the shock absorbers.'']'
$ws.Cells.Item(21, 1).Value = 'Verse 20'
$ws.Cells.Item(21, 2).Value = '["He ain''t crude, look at that, look at that"]'
$ws.Cells.Item(22, 1).Value = 'Verse 21'
$ws.Cells.Item(22, 2).Value = '["He ain''t lewd, look at that, look at that"]'
$ws.Cells.Item(23, 1).Value = 'Verse 22'
$ws.Cells.Item(23, 2).Value = '[''(Look at that, look at that)'']'
$ws.Cells.Item(24, 1).Value = 'Verse 23'
$ws.Cells.Item(24, 2).Value = '[''(Look at that, look at that)'']'
$ws.Cells.Item(25, 1).Value = 'Verse 24'
$ws.Cells.Item(25, 2).Value = '["He''s always makin'' the news"]'
$ws.Cells.Item(26, 1).Value = 'Verse 25'
$ws.Cells.Item(26, 2).Value = '["Wearin'' just his tennis shoes"]'
$ws.Cells.Item(27, 1).Value = 'Verse 26'
$ws.Cells.Item(27, 2).Value = '[''Guess you could call him unique'']'
$ws.Cells.Item(28, 1).Value = 'Verse 27'
$ws.Cells.Item(28, 2).Value = '[''Yeah, I did.'']'
$ws.Cells.Item(29, 1).Value = 'Verse 28'
$ws.Cells.Item(29, 2).Value = '["Half time, I''s just goin'' down thar to get Ethel a snow cone."]'
$ws.Cells.Item(30, 1).Value = 'Verse 29'
$ws.Cells.Item(30, 2).Value = '[''And here he come, right out of the cheap seats, dribbling, right down the middle of the court.'']'
$ws.Cells.Item(31, 1).Value = 'Verse 30'
$ws.Cells.Item(31, 2).Value = '["Didn''t have on nothing but his PF''s."]'
$ws.Cells.Item(32, 1).Value = 'Verse 31'
$ws.Cells.Item(32, 2).Value = '[''Made a hook shot and got out through the concessions stand.'']'
$ws.Cells.Item(33, 1).Value = 'Verse 32'
$ws.Cells.Item(33, 2).Value = '[''But it was too late.'']'
$ws.Cells.Item(34, 1).Value = 'Verse 33'
$ws.Cells.Item(34, 2).Value = '["She''d already got a free shot."]'
$ws.Cells.Item(35, 1).Value = 'Verse 34'
$ws.Cells.Item(35, 2).Value = '[''(Look at that, look at that)'']'
$ws.Cells.Item(36, 1).Value = 'Verse 35'
$ws.Cells.Item(36, 2).Value = '[''Of his anatomy'']'
$ws.Cells.Item(37, 1).Value = 'Verse 36'
$ws.Cells.Item(37, 2).Value = '["He''s gonna give us a peek"]'
$ws.Cells.Item(38, 1).Value = 'Verse 37'
$ws.Cells.Item(38, 2).Value = '[''Here he comes again.'']'
$ws.Cells.Item(39, 1).Value = 'Verse 38'
$ws.Cells.Item(39, 2).Value = '["Who''s that with him?"]'
$ws.Cells.Item(40, 1).Value = 'Verse 39'
$ws.Cells.Item(40, 2).Value = '[''Ethel?'']'
$ws.Cells.Item(41, 1).Value = 'Verse 40'
$ws.Cells.Item(41, 2).Value = '[''Is that you, Ethel?'']'
$ws.Cells.Item(42, 1).Value = 'Verse 41'
$ws.Cells.Item(42, 2).Value = '["What do you think you''re"]'
$ws.Cells.Item(43, 1).Value = 'Verse 42'
$ws.Cells.Item(43, 2).Value = '["doin''?"]'
$ws.Cells.Item(44, 1).Value = 'Verse 43'
$ws.Cells.Item(44, 2).Value = '[''You git your'']'
$ws.Cells.Item(45, 1).Value = 'Verse 44'
$ws.Cells.Item(45, 2).Value = '[]'
$ws.Cells.Item(46, 1).Value = 'Verse 45'
$ws.Cells.Item(46, 2).Value = '["Where you goin''?"]'
$ws.Cells.Item(47, 1).Value = 'Verse 46'
$ws.Cells.Item(47, 2).Value = '[''Ethel, you shameless'']'
$ws.Cells.Item(48, 1).Value = 'Verse 47'
$ws.Cells.Item(48, 2).Value = '[''hussy!'']'
$ws.Cells.Item(49, 1).Value = 'Verse 48'
$ws.Cells.Item(49, 2).Value = '[''Ethel!'']'
$ws.Cells.Item(50, 1).Value = 'Verse 49'
$ws.Cells.Item(50, 2).Value = '[''Ethelllllll!!'']'
$ws.Cells.Item(51, 1).Value = 'Verse 50'
$ws.Cells.Item(51, 2).Value = '[''!'']'
$ws.Cells.Item(52, 1).Value = 'Pre-or-Post-Chorus 1'
$ws.Cells.Item(52, 2).Value = 'Pardon me, sir, did you see what happened?'
$ws.Cells.Item(53, 1).Value = 'Pre-or-Post-Chorus 2'
$ws.Cells.Item(53, 2).Value = 'Yeah, I did.'
$ws.Cells.Item(54, 1).Value = 'Pre-or-Post-Chorus 3'
$ws.Cells.Item(54, 2).Value = 'He''s just as proud as he can be'
$ws.Cells.Item(55, 1).Value = 'Pre-or-Post-Chorus 4'
$ws.Cells.Item(55, 2).Value = 'Of his anatomy'
$ws.Cells.Item(56, 1).Value = 'Pre-or-Post-Chorus 5'
$ws.Cells.Item(56, 2).Value = 'He likes to show off his physique'
$ws.Cells.Item(57, 1).Value = 'Pre-or-Post-Chorus 6'
$ws.Cells.Item(57, 2).Value = 'If there''s an audience to be found'
$ws.Cells.Item(58, 1).Value = 'Pre-or-Post-Chorus 7'
$ws.Cells.Item(58, 2).Value = 'He''ll be streakin'' around'
$ws.Cells.Item(59, 1).Value = 'Pre-or-Post-Chorus 8'
$ws.Cells.Item(59, 2).Value = 'Invitin'' public critique'
$ws.Cells.Item(60, 1).Value = 'Pre-or-Post-Chorus 9'
$ws.Cells.Item(60, 2).Value = 'But it was too late.'
$ws.Cells.Item(61, 1).Value = 'Pre-or-Post-Chorus 10'
$ws.Cells.Item(61, 2).Value = 'Ethel!'
$ws.Cells.Item(62, 1).Value = 'Chorus 1'
$ws.Cells.Item(62, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(63, 1).Value = 'Chorus 2'
$ws.Cells.Item(63, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(64, 1).Value = 'Chorus 3'
$ws.Cells.Item(64, 2).Value = 'Oh, yes, they call him the Streak'
$ws.Cells.Item(65, 1).Value = 'Chorus 4'
$ws.Cells.Item(65, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(66, 1).Value = 'Chorus 5'
$ws.Cells.Item(66, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(67, 1).Value = 'Chorus 6'
$ws.Cells.Item(67, 2).Value = 'Oh, yes, they call him the Streak'
$ws.Cells.Item(68, 1).Value = 'Chorus 7'
$ws.Cells.Item(68, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(69, 1).Value = 'Chorus 8'
$ws.Cells.Item(69, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(70, 1).Value = 'Chorus 9'
$ws.Cells.Item(70, 2).Value = 'Oh, yes, they call him the Streak'
$ws.Cells.Item(71, 1).Value = 'Chorus 10'
$ws.Cells.Item(71, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(72, 1).Value = 'Chorus 11'
$ws.Cells.Item(72, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(73, 1).Value = 'Chorus 12'
$ws.Cells.Item(73, 2).Value = 'Oh, yes, they call him the Streak'
$ws.Cells.Item(74, 1).Value = 'Chorus 13'
$ws.Cells.Item(74, 2).Value = '(Look at that, look at that)'
$ws.Cells.Item(75, 1).Value = 'Chorus 14'
$ws.Cells.Item(75, 2).Value = '(Look at that, look at that)'
